$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-64 down to 43-65.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record.
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44518
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 100112022
$ws.Range("G42").Value = "Arveja Verde"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 600
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = 15000
$ws.Range("N42").Value = "$/saco 25 kilos"
$ws.Range("O42").Value = "Región del Maule"
$ws.Range("P42").Value = 600
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"

# Make sure the date cell uses the same number format style as the rest of column D.
$ws.Range("D42").NumberFormat = $ws.Range("D43").NumberFormat
